# Apply hybrid bold + color (2C3E50) highlighting to quantitative impact
# metrics (percentages, dollar amounts) across the resume body, matching
# the commit's "quantitative metrics highlighting" feature.
#
# Each entry below identifies the 1-based paragraph index in the document
# and the literal metric text to find and highlight within that paragraph.
# wdColor values are BGR-packed integers; 0x503E2C == RGB(2C,3E,50).

$d = $word.ActiveDocument
$metricColor = 5258796  # 0x503E2C -> w:color val="2C3E50"

$targets = @(
    @{p = 10; t = '23%'},
    @{p = 10; t = '64%'},

    @{p = 12; t = '±4.2%'},
    @{p = 12; t = '±2.1%'},
    @{p = 12; t = '71%'},
    @{p = 12; t = '87%'},

    @{p = 13; t = '73.5%'},
    @{p = 13; t = '$4.7M'},

    @{p = 14; t = '$2'},

    @{p = 34; t = '57%'},

    @{p = 50; t = '±4.2%'},
    @{p = 50; t = '±2.1%'},

    @{p = 51; t = '71%'},
    @{p = 51; t = '87%'},

    @{p = 52; t = '34%'},
    @{p = 52; t = '28%'}
)

foreach ($item in $targets) {
    $para = $d.Paragraphs.Item($item.p)
    $rng = $para.Range
    $found = $rng.Find.Execute($item.t, $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
    if ($found) {
        $rng.Font.Bold = 1
        $rng.Font.Color = $metricColor
    }
}
